$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, copying the header style (bold font + border + alignment)
# from the existing last header cell (AB1) so they match the other headers.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Team record columns for every data row (2-38): Wins, Losses, Ties.
$lastRow = 38
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 84
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 0
}
